# Insert a new data row at row 117 (pushing existing rows 117-172 down to
# 118-173), then populate the new row with the additional weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 6
$ws.Cells.Item(117, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(117, 3).Value = "Metropolitana"
$ws.Cells.Item(117, 4).Value = 44523
$ws.Cells.Item(117, 5).Value = 13
$ws.Cells.Item(117, 6).Value = 100112022
$ws.Cells.Item(117, 7).Value = "Arveja Verde"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 200
$ws.Cells.Item(117, 11).Value = 14000
$ws.Cells.Item(117, 12).Value = 15000
$ws.Cells.Item(117, 13).Value = 14600
$ws.Cells.Item(117, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(117, 15).Value = "Región del Maule"
$ws.Cells.Item(117, 16).Value = 584
$ws.Cells.Item(117, 17).Value = 25
$ws.Cells.Item(117, 18).Value = "Hortaliza"
